$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new test case row (row 15): "Create Project Object" ---
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "test_<CreateProject>"
$ws.Range("C15").Value = "This is to test whether a project can be created successfully"
$ws.Range("D15").Value = "title: ""New Project""`r`ndescription: ""Testing project""`r`ntechnology: ""Project""`r`nimage: ""img/project1.png"""
$ws.Range("D15").WrapText = $true
$ws.Range("E15").Value = "A project object is created"

# Match the row height Excel computed for the wrapped 4-line entry
$ws.Rows.Item(15).RowHeight = 58.3

# --- Update the active selection / view to reflect where the author ended up ---
$ws.Range("G18").Select()
